# Add data for 2021-12-04
# - Rename sheet / update header label from "...November 25" to "...November 26"
# - Update the carjacking-by-neighborhood-by-month counts for the
#   "November 2021 (through November 26)" column (column B) and a handful
#   of historical-month columns that received revised counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab and update the running-total header text.
$ws.Name = "Through 2021-11-26"
$ws.Range("B1").Value = "November 2021 (through November 26)"

# Updated / newly-populated cell counts.
$ws.Range("E3").Value = 14
$ws.Range("X3").Value = 6
$ws.Range("AI3").Value = 6
$ws.Range("M4").Value = 12
$ws.Range("BE4").Value = 4
$ws.Range("B5").Value = 5
$ws.Range("E5").Value = 9
$ws.Range("AT5").Value = 4
$ws.Range("AT7").Value = 4
$ws.Range("BP8").Value = 1
$ws.Range("B10").Value = 2
$ws.Range("AI13").Value = 3
$ws.Range("B17").Value = 4
$ws.Range("B20").Value = 5
$ws.Range("AI21").Value = 1
$ws.Range("BP22").Value = 1
$ws.Range("B25").Value = 2
$ws.Range("BE38").Value = 1
$ws.Range("AI44").Value = 1
$ws.Range("BE48").Value = 4
$ws.Range("BP48").Value = 3
$ws.Range("AI49").Value = 1
$ws.Range("AT52").Value = 2
$ws.Range("AI59").Value = 2
$ws.Range("M64").Value = 3
$ws.Range("AI88").Value = 1
$ws.Range("BE99").Value = 1

Write-Output "Applied 2021-12-04 data update"
